# Add the "2022-Q4" quarterly report sheet to the 600704-物产中大 holdings
# workbook ("feat: add 2022-Q4 data").
#
# What changes:
#   1. A brand new worksheet named "2022-Q4" is inserted right after "总计"
#      (and before "2022-Q3"), containing the fund-holding detail rows for
#      that quarter.
#   2. The "总计" (summary) sheet gets a new row inserted right under its
#      header for "2022-Q4", every existing row shifts down by one, and the
#      running index in column A is renumbered 0..n accordingly.
#   3. The previously-active tab ("2020-Q4") stays the active/selected tab.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 1. Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Add($null, $totalSheet)
$q4Sheet.Name = "2022-Q4"

# Match the page-margin convention used by the other quarterly detail sheets
# (0.75in / 1in / 0.5in margins, expressed in points: 1in = 72pt).
$q4Sheet.PageSetup.LeftMargin = 54
$q4Sheet.PageSetup.RightMargin = 54
$q4Sheet.PageSetup.TopMargin = 72
$q4Sheet.PageSetup.BottomMargin = 72
$q4Sheet.PageSetup.HeaderMargin = 36
$q4Sheet.PageSetup.FooterMargin = 36

# Reuse the bordered/bold header style and the bordered index-column style
# that the rest of the workbook already uses (both happen to be the same
# style, carried by 总计!B1 and 总计!A2) so the new sheet matches visually.
$sharedStyleCell = $totalSheet.Cells.Item(1, 2)
$headerRng = $q4Sheet.Range("B1:H1")
$sharedStyleCell.Copy()
$headerRng.PasteSpecial($xlPasteFormats)

$idxRng = $q4Sheet.Range("A2:A4")
$sharedStyleCell.Copy()
$idxRng.PasteSpecial($xlPasteFormats)
$q4Sheet.Application.CutCopyMode = $false

# Header row.
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $q4Sheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}

# Columns B-G (fund code, name, scale, position, ratio, market value) are
# stored as plain text in every quarterly sheet in this workbook (so
# leading zeros in fund codes and trailing zeros like "3.70"/"98.00"
# survive) - only the running index (A) and the rank (H) are numeric.
# (Bounded range, not a whole-column reference, so the text format is
# actually honoured when the values are written below.)
$q4Sheet.Range("B2:G4").NumberFormat = "@"

# Detail rows: code, name, scale, total stock position, position ratio,
# held market value (yi), position rank.
$q4Data = @(
    @("515760", "华夏中证浙江国资创新发展ETF", "2.14", "99.05", "3.70", "0.0792", 8),
    @("516910", "富国中证现代物流ETF",         "1.09", "98.85", "3.77", "0.0411", 9),
    @("516530", "银华中证现代物流ETF",         "0.85", "98.00", "3.72", "0.0316", 9)
)

for ($i = 0; $i -lt $q4Data.Length; $i++) {
    $r = $i + 2
    $entry = $q4Data[$i]
    $q4Sheet.Cells.Item($r, 1).Value = $i
    $q4Sheet.Cells.Item($r, 2).Value = $entry[0]
    $q4Sheet.Cells.Item($r, 3).Value = $entry[1]
    $q4Sheet.Cells.Item($r, 4).Value = $entry[2]
    $q4Sheet.Cells.Item($r, 5).Value = $entry[3]
    $q4Sheet.Cells.Item($r, 6).Value = $entry[4]
    $q4Sheet.Cells.Item($r, 7).Value = $entry[5]
    $q4Sheet.Cells.Item($r, 8).Value = $entry[6]
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" summary table: add the 2022-Q4 row, shift the rest
#    down and renumber the running index.
# ---------------------------------------------------------------------------
$summaryDates  = @("2022-Q4","2022-Q3","2022-Q2","2022-Q1","2021-Q4","2021-Q3","2021-Q2","2021-Q1","2020-Q4")
$summaryCounts = @(3, 7, 11, 8, 2, 7, 11, 5, 5)
$summaryValues = @(0.15, 0.17, 0.43, 0.42, 0.04, 4.07, 4.84, 0.13, 0.16)

for ($i = 0; $i -lt $summaryDates.Length; $i++) {
    $r = $i + 2
    $totalSheet.Cells.Item($r, 1).Value = $i
    $totalSheet.Cells.Item($r, 2).Value = $summaryDates[$i]
    $totalSheet.Cells.Item($r, 3).Value = $summaryCounts[$i]
    $totalSheet.Cells.Item($r, 4).Value = $summaryValues[$i]
}

# The newly-written last row (A10) has no pre-existing style to inherit, so
# copy the running-index formatting from the row above it.
$lastRow = $summaryDates.Length + 1
$totalSheet.Cells.Item($lastRow - 1, 1).Copy($totalSheet.Cells.Item($lastRow, 1))
$totalSheet.Cells.Item($lastRow, 1).Value = $summaryDates.Length - 1

# ---------------------------------------------------------------------------
# 3. Restore "2020-Q4" as the active/selected tab (it was the tab open when
#    the workbook was last saved, and stays so after this edit).
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2020-Q4").Activate()
